# Apply trade #55 close update across the workbook.
$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = -0.91
$summary.Range("B6").Value = 55
$summary.Range("B9").Value = 40

# ---- Strategy Status sheet ----
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 55
$status.Range("G4").Value = 40

# ---- New trade row data (same for "All Trades" and "MarketMaking" sheets) ----
$newRow = @(55, "2026-02-17", "13:29:33", "MarketMaking", "UP", 0.97, 0.97, "CLOSED", 0, 0, 97.5, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowIndex = 56
    # Force the Date/Time columns (B, C) to be stored as plain text so Excel
    # does not auto-convert them into date/time serial numbers.
    $ws.Range("B$rowIndex").NumberFormat = "@"
    $ws.Range("C$rowIndex").NumberFormat = "@"
    for ($col = 1; $col -le $newRow.Length; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $newRow[$col - 1]
    }
}
